$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = [double]"1.724404365356792e-05"
$ws.Range("D2").Value = [double]"9.871345059498339e-07"
$ws.Range("E2").Value = [double]"7.033013450984647e-07"

$ws.Range("B3").Value = [double]"1.724404365356792e-05"
$ws.Range("D3").Value = [double]"1.974535298499556e-06"
$ws.Range("E3").Value = [double]"2.213103850403893"

$ws.Range("B4").Value = [double]"9.871345059498339e-07"
$ws.Range("C4").Value = [double]"1.974535298499556e-06"
$ws.Range("E4").Value = [double]"0.8572568171758366"

$ws.Range("B5").Value = [double]"7.033013450984647e-07"
$ws.Range("C5").Value = [double]"2.213103850403893"
$ws.Range("D5").Value = [double]"0.8572568171758366"
